$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 156.82608
$ws.Cells.Item(33, 9).Value = 149.8421
$ws.Cells.Item(33, 11).Value = 149.8421
$ws.Cells.Item(33, 13).Value = 79.15790000000001
$ws.Cells.Item(34, 8).Value = 1872
$ws.Cells.Item(34, 9).Value = 1872
$ws.Cells.Item(34, 11).Value = 1872
$ws.Cells.Item(34, 13).Value = -1669
$ws.Cells.Item(36, 8).Value = 1872
$ws.Cells.Item(36, 9).Value = 1872
$ws.Cells.Item(36, 11).Value = 1872
$ws.Cells.Item(36, 13).Value = -1157
$ws.Cells.Item(98, 8).Value = 1436.7778
$ws.Cells.Item(98, 9).Value = 1324.1333
$ws.Cells.Item(98, 10).Value = 2000
$ws.Cells.Item(98, 11).Value = 1324.1333
$ws.Cells.Item(98, 12).Value = 2000
$ws.Cells.Item(98, 13).Value = 173.8667
$ws.Cells.Item(98, 14).Value = -4996
$ws.Cells.Item(99, 8).Value = 2621.5386
$ws.Cells.Item(99, 9).Value = 214.66667
$ws.Cells.Item(99, 10).Value = 3343.6
$ws.Cells.Item(99, 11).Value = 644.00001
$ws.Cells.Item(99, 12).Value = 10030.8
$ws.Cells.Item(99, 13).Value = 853.99999
$ws.Cells.Item(99, 14).Value = -13026.8
$ws.Cells.Item(122, 8).Value = 1436.7778
$ws.Cells.Item(122, 9).Value = 1324.1333
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 3972.3999
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -1522.3999
$ws.Cells.Item(122, 14).Value = -10900
$ws.Cells.Item(124, 8).Value = 72300
$ws.Cells.Item(124, 10).Value = 72300
$ws.Cells.Item(124, 12).Value = 72300
$ws.Cells.Item(124, 14).Value = -82120
$ws.Cells.Item(126, 8).Value = 46000
$ws.Cells.Item(126, 10).Value = 46000
$ws.Cells.Item(126, 12).Value = 46000
$ws.Cells.Item(126, 14).Value = -55880
$ws.Cells.Item(130, 8).Value = 38997.5
$ws.Cells.Item(130, 10).Value = 38997.5
$ws.Cells.Item(130, 12).Value = 38997.5
$ws.Cells.Item(130, 14).Value = -49037.5
$ws.Cells.Item(137, 8).Value = 1527
$ws.Cells.Item(137, 9).Value = 1958.6
$ws.Cells.Item(137, 10).Value = 1354.36
$ws.Cells.Item(137, 11).Value = 5875.799999999999
$ws.Cells.Item(137, 12).Value = 4063.08
$ws.Cells.Item(137, 13).Value = -3325.799999999999
$ws.Cells.Item(137, 14).Value = -9163.08
$ws.Cells.Item(141, 8).Value = 5199
$ws.Cells.Item(141, 9).Value = 2309.875
$ws.Cells.Item(141, 10).Value = 8088.125
$ws.Cells.Item(141, 11).Value = 6929.625
$ws.Cells.Item(141, 12).Value = 24264.375
$ws.Cells.Item(141, 13).Value = -1749.625
$ws.Cells.Item(141, 14).Value = -34624.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 12197754
$ws.Cells.Item(74, 9).Value = 1179.2
$ws.Cells.Item(74, 10).Value = 23813540
$ws.Cells.Item(74, 11).Value = 1179.2
$ws.Cells.Item(74, 12).Value = 23813540
$ws.Cells.Item(74, 13).Value = -305.2
$ws.Cells.Item(74, 14).Value = -23815288
$ws.Cells.Item(77, 8).Value = 12197754
$ws.Cells.Item(77, 9).Value = 1179.2
$ws.Cells.Item(77, 10).Value = 23813540
$ws.Cells.Item(77, 11).Value = 5896
$ws.Cells.Item(77, 12).Value = 119067700
$ws.Cells.Item(77, 13).Value = -1528
$ws.Cells.Item(77, 14).Value = -119076436
$ws.Cells.Item(122, 8).Value = 42483.56
$ws.Cells.Item(122, 9).Value = 73479.64
$ws.Cells.Item(122, 10).Value = 3034
$ws.Cells.Item(122, 11).Value = 220438.92
$ws.Cells.Item(122, 12).Value = 9102
$ws.Cells.Item(122, 13).Value = -217988.92
$ws.Cells.Item(122, 14).Value = -14002
$ws.Cells.Item(132, 8).Value = 1285915.9
$ws.Cells.Item(132, 9).Value = 3310.319
$ws.Cells.Item(132, 11).Value = 9930.957
$ws.Cells.Item(132, 13).Value = -7400.957

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(126, 8).Value = 60780
$ws.Cells.Item(126, 10).Value = 60780
$ws.Cells.Item(126, 12).Value = 60780
$ws.Cells.Item(126, 14).Value = -70660

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1757.7778
$ws.Cells.Item(16, 9).Value = 2384.4
$ws.Cells.Item(16, 10).Value = 974.5
$ws.Cells.Item(16, 11).Value = 2384.4
$ws.Cells.Item(16, 12).Value = 974.5
$ws.Cells.Item(16, 13).Value = -2097.4
$ws.Cells.Item(16, 14).Value = -1548.5
$ws.Cells.Item(31, 8).Value = 7095.479
$ws.Cells.Item(31, 9).Value = 1799.3
$ws.Cells.Item(31, 10).Value = 8489.210999999999
$ws.Cells.Item(31, 11).Value = 1799.3
$ws.Cells.Item(31, 12).Value = 8489.210999999999
$ws.Cells.Item(31, 13).Value = -1504.3
$ws.Cells.Item(31, 14).Value = -9079.210999999999
$ws.Cells.Item(32, 8).Value = 1650
$ws.Cells.Item(32, 9).Value = 1000
$ws.Cells.Item(32, 10).Value = 2300
$ws.Cells.Item(32, 11).Value = 1000
$ws.Cells.Item(32, 12).Value = 2300
$ws.Cells.Item(32, 13).Value = -684
$ws.Cells.Item(32, 14).Value = -2932
$ws.Cells.Item(34, 8).Value = 7095.479
$ws.Cells.Item(34, 9).Value = 1799.3
$ws.Cells.Item(34, 10).Value = 8489.210999999999
$ws.Cells.Item(34, 11).Value = 1799.3
$ws.Cells.Item(34, 12).Value = 8489.210999999999
$ws.Cells.Item(34, 13).Value = -1597.3
$ws.Cells.Item(34, 14).Value = -8893.210999999999
$ws.Cells.Item(113, 8).Value = 1757.7778
$ws.Cells.Item(113, 9).Value = 2384.4
$ws.Cells.Item(113, 10).Value = 974.5
$ws.Cells.Item(113, 11).Value = 2384.4
$ws.Cells.Item(113, 12).Value = 974.5
$ws.Cells.Item(113, 13).Value = -214.4000000000001
$ws.Cells.Item(113, 14).Value = -5314.5
$ws.Cells.Item(132, 8).Value = 2958.3125
$ws.Cells.Item(132, 9).Value = 2775.9092
$ws.Cells.Item(132, 10).Value = 3359.6
$ws.Cells.Item(132, 11).Value = 8327.7276
$ws.Cells.Item(132, 12).Value = 10078.8
$ws.Cells.Item(132, 13).Value = -5797.7276
$ws.Cells.Item(132, 14).Value = -15138.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 982.8570999999999
$ws.Cells.Item(22, 9).Value = 600
$ws.Cells.Item(22, 10).Value = 1028.8
$ws.Cells.Item(22, 11).Value = 1800
$ws.Cells.Item(22, 12).Value = 3086.4
$ws.Cells.Item(22, 13).Value = -1631
$ws.Cells.Item(22, 14).Value = -3424.4
$ws.Cells.Item(27, 8).Value = 982.8570999999999
$ws.Cells.Item(27, 9).Value = 600
$ws.Cells.Item(27, 10).Value = 1028.8
$ws.Cells.Item(27, 11).Value = 1800
$ws.Cells.Item(27, 12).Value = 3086.4
$ws.Cells.Item(27, 13).Value = -1698
$ws.Cells.Item(27, 14).Value = -3290.4
$ws.Cells.Item(36, 8).Value = 2176.923
$ws.Cells.Item(36, 9).Value = 1433.3334
$ws.Cells.Item(36, 10).Value = 2400
$ws.Cells.Item(36, 11).Value = 4300.0002
$ws.Cells.Item(36, 12).Value = 7200
$ws.Cells.Item(36, 13).Value = -4131.0002
$ws.Cells.Item(36, 14).Value = -7538
$ws.Cells.Item(131, 8).Value = 5436.6665
$ws.Cells.Item(131, 10).Value = 6003.7036
$ws.Cells.Item(131, 12).Value = 18011.1108
$ws.Cells.Item(131, 14).Value = -28091.1108
$ws.Cells.Item(141, 8).Value = 4140.8184
$ws.Cells.Item(141, 9).Value = 2026.5834
$ws.Cells.Item(141, 10).Value = 5348.952
$ws.Cells.Item(141, 11).Value = 6079.7502
$ws.Cells.Item(141, 12).Value = 16046.856
$ws.Cells.Item(141, 13).Value = -899.7502000000004
$ws.Cells.Item(141, 14).Value = -26406.856

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 76670.664
$ws.Cells.Item(12, 10).Value = 76670.664
$ws.Cells.Item(12, 12).Value = 76670.664
$ws.Cells.Item(12, 14).Value = -76950.664
$ws.Cells.Item(62, 8).Value = 29996
$ws.Cells.Item(62, 10).Value = 29996
$ws.Cells.Item(62, 12).Value = 29996
$ws.Cells.Item(62, 14).Value = -31368
$ws.Cells.Item(65, 8).Value = 29996
$ws.Cells.Item(65, 10).Value = 29996
$ws.Cells.Item(65, 12).Value = 89988
$ws.Cells.Item(65, 14).Value = -96852
$ws.Cells.Item(123, 8).Value = 8318.385
$ws.Cells.Item(123, 10).Value = 8318.385
$ws.Cells.Item(123, 12).Value = 8318.385
$ws.Cells.Item(123, 14).Value = -13218.385
$ws.Cells.Item(132, 8).Value = 4294.5
$ws.Cells.Item(132, 9).Value = 3561.1
$ws.Cells.Item(132, 10).Value = 6128
$ws.Cells.Item(132, 11).Value = 10683.3
$ws.Cells.Item(132, 12).Value = 18384
$ws.Cells.Item(132, 13).Value = -8153.299999999999
$ws.Cells.Item(132, 14).Value = -23444

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 2998.2
$ws.Cells.Item(122, 9).Value = 2312.5715
$ws.Cells.Item(122, 10).Value = 3598.125
$ws.Cells.Item(122, 11).Value = 6937.7145
$ws.Cells.Item(122, 12).Value = 10794.375
$ws.Cells.Item(122, 13).Value = -4487.7145
$ws.Cells.Item(122, 14).Value = -15694.375
$ws.Cells.Item(137, 8).Value = 26777.777
$ws.Cells.Item(137, 10).Value = 27250
$ws.Cells.Item(137, 12).Value = 27250
$ws.Cells.Item(137, 14).Value = -37450

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(47, 8).Value = 15000
$ws.Cells.Item(47, 10).Value = 15000
$ws.Cells.Item(47, 12).Value = 15000
$ws.Cells.Item(47, 14).Value = -16144
$ws.Cells.Item(123, 8).Value = 40429
$ws.Cells.Item(123, 10).Value = 40429
$ws.Cells.Item(123, 12).Value = 40429
$ws.Cells.Item(123, 14).Value = -50229
$ws.Cells.Item(136, 8).Value = 5752.5835
$ws.Cells.Item(136, 9).Value = 6608.6665
$ws.Cells.Item(136, 10).Value = 4896.5
$ws.Cells.Item(136, 11).Value = 19825.9995
$ws.Cells.Item(136, 12).Value = 14689.5
$ws.Cells.Item(136, 13).Value = -17275.9995
$ws.Cells.Item(136, 14).Value = -19789.5
